# Update the "Bico" and "Tanque" sheets:
#  - Rename the "Obs" header to "Obs_relatorio"
#  - Insert a new trailing "Obs_sped" column (left blank on every data row)
#  - Replace every "VERDADEIRO" validation flag with the full success message

$wb = $excel.ActiveWorkbook

$message = "Validado com sucesso! Nenhuma divergência entre o SPED e o relatório foi encontrada!"

# --- Sheet "Bico" (H = Obs column, I = new Obs_sped column, data rows 2-13) ---
$wsBico = $wb.Worksheets.Item("Bico")

$wsBico.Range("H1").Value = "Obs_relatorio"
$wsBico.Range("I1").Value = "Obs_sped"

for ($r = 2; $r -le 13; $r++) {
    $wsBico.Cells.Item($r, 8).Value = $message
    # A lone "'" gives an empty, text-typed cell (like typing ' + Enter in Excel)
    # instead of leaving the cell completely blank/unset.
    $wsBico.Cells.Item($r, 9).Value = "'"
    $wsBico.Cells.Item($r, 9).ClearFormats()
}

# --- Sheet "Tanque" (F = Obs column, G = new Obs_sped column, data rows 2-9) ---
$wsTanque = $wb.Worksheets.Item("Tanque")

$wsTanque.Range("F1").Value = "Obs_relatorio"
$wsTanque.Range("G1").Value = "Obs_sped"

for ($r = 2; $r -le 9; $r++) {
    $wsTanque.Cells.Item($r, 6).Value = $message
    $wsTanque.Cells.Item($r, 7).Value = "'"
    $wsTanque.Cells.Item($r, 7).ClearFormats()
}
